$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 2 (JOBY) updates
$ws.Range("D2").Value = 15.6
$ws.Range("E2").Value = 56.4
$ws.Range("F2").Value = 10.45
$ws.Range("H2").Value = 53
$ws.Range("J2").Value = 66
$ws.Range("N2").Value = 53.71147335634279

# Row 3 (ACHR) updates
$ws.Range("D3").Value = 8.68
$ws.Range("E3").Value = 58.8
$ws.Range("F3").Value = 15.82
$ws.Range("N3").Value = 53.71147335634279
